$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 445.14285
$ws.Range("I2").Value = 399.33334
$ws.Range("K2").Value = 399.33334
$ws.Range("M2").Value = -286.33334
$ws.Range("H4").Value = 203.33333
$ws.Range("I4").Value = 203.33333
$ws.Range("K4").Value = 203.33333
$ws.Range("M4").Value = -89.33332999999999
$ws.Range("H17").Value = 2572269.5
$ws.Range("J17").Value = 2572269.5
$ws.Range("L17").Value = 7716808.5
$ws.Range("N17").Value = -7717144.5
$ws.Range("H19").Value = 8296.333
$ws.Range("I19").Value = 1889.5
$ws.Range("J19").Value = 11499.75
$ws.Range("K19").Value = 1889.5
$ws.Range("L19").Value = 11499.75
$ws.Range("M19").Value = -1714.5
$ws.Range("N19").Value = -11849.75
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = $null
$ws.Range("N32").Value = $null
$ws.Range("H33").Value = 226.18182
$ws.Range("I33").Value = 177.375
$ws.Range("J33").Value = 356.33334
$ws.Range("K33").Value = 177.375
$ws.Range("L33").Value = 356.33334
$ws.Range("M33").Value = 51.625
$ws.Range("N33").Value = -814.33334
$ws.Range("H40").Value = 5788.8
$ws.Range("J40").Value = 4174.625
$ws.Range("L40").Value = 4174.625
$ws.Range("N40").Value = -4524.625
$ws.Range("H43").Value = 4284.5
$ws.Range("J43").Value = 5999.5
$ws.Range("L43").Value = 5999.5
$ws.Range("N43").Value = -6137.5
$ws.Range("H51").Value = 8917.0625
$ws.Range("J51").Value = 9266.167
$ws.Range("L51").Value = 9266.167
$ws.Range("N51").Value = -10234.167
$ws.Range("H58").Value = 1366.5883
$ws.Range("I58").Value = 475.63635
$ws.Range("K58").Value = 1426.90905
$ws.Range("M58").Value = -1276.90905
$ws.Range("H61").Value = 6832.533
$ws.Range("I61").Value = 7142
$ws.Range("K61").Value = 21426
$ws.Range("M61").Value = -21254
$ws.Range("H62").Value = 13897379
$ws.Range("I62").Value = 33341132
$ws.Range("K62").Value = 33341132
$ws.Range("M62").Value = -33340508
$ws.Range("H64").Value = 3807.1538
$ws.Range("I64").Value = 3249.4167
$ws.Range("J64").Value = 10500
$ws.Range("K64").Value = 3249.4167
$ws.Range("L64").Value = 10500
$ws.Range("M64").Value = -3001.4167
$ws.Range("N64").Value = -10996
$ws.Range("H65").Value = 13897379
$ws.Range("I65").Value = 33341132
$ws.Range("K65").Value = 166705660
$ws.Range("M65").Value = -166702540
$ws.Range("H67").Value = 3807.1538
$ws.Range("I67").Value = 3249.4167
$ws.Range("J67").Value = 10500
$ws.Range("K67").Value = 3249.4167
$ws.Range("L67").Value = 10500
$ws.Range("M67").Value = -2391.4167
$ws.Range("N67").Value = -12216
$ws.Range("H74").Value = 5571
$ws.Range("I74").Value = 5571
$ws.Range("K74").Value = 5571
$ws.Range("M74").Value = -4635
$ws.Range("H77").Value = 5571
$ws.Range("I77").Value = 5571
$ws.Range("K77").Value = 27855
$ws.Range("M77").Value = -23175
$ws.Range("H86").Value = 6366.3335
$ws.Range("I86").Value = 5050
$ws.Range("J86").Value = 8999
$ws.Range("K86").Value = 5050
$ws.Range("L86").Value = 8999
$ws.Range("M86").Value = -3927
$ws.Range("N86").Value = -11245
$ws.Range("H88").Value = 3999
$ws.Range("J88").Value = 3999
$ws.Range("L88").Value = 3999
$ws.Range("N88").Value = -4811
$ws.Range("H89").Value = 6366.3335
$ws.Range("I89").Value = 5050
$ws.Range("J89").Value = 8999
$ws.Range("K89").Value = 25250
$ws.Range("L89").Value = 44995
$ws.Range("M89").Value = -19634
$ws.Range("N89").Value = -56227
$ws.Range("H91").Value = 3999
$ws.Range("J91").Value = 3999
$ws.Range("L91").Value = 3999
$ws.Range("N91").Value = -6807
$ws.Range("H96").Value = 646.2778
$ws.Range("I96").Value = 737.4667
$ws.Range("J96").Value = 190.33333
$ws.Range("K96").Value = 2212.4001
$ws.Range("L96").Value = 570.99999
$ws.Range("M96").Value = -839.4000999999998
$ws.Range("N96").Value = -3316.99999
$ws.Range("H100").Value = 7695.8335
$ws.Range("I100").Value = 3235
$ws.Range("J100").Value = 30000
$ws.Range("K100").Value = 3235
$ws.Range("L100").Value = 30000
$ws.Range("M100").Value = -2694
$ws.Range("N100").Value = -31082
$ws.Range("H103").Value = 561.44446
$ws.Range("J103").Value = 764.375
$ws.Range("L103").Value = 2293.125
$ws.Range("N103").Value = -3465.125
$ws.Range("H112").Value = 3874365.5
$ws.Range("J112").Value = 3874365.5
$ws.Range("L112").Value = 11623096.5
$ws.Range("N112").Value = -11625312.5
$ws.Range("H113").Value = 9708.223
$ws.Range("J113").Value = 4567.1665
$ws.Range("L113").Value = 4567.1665
$ws.Range("N113").Value = -11075.1665
$ws.Range("H132").Value = 3580.4546
$ws.Range("I132").Value = 2826.2415
$ws.Range("K132").Value = 8478.7245
$ws.Range("M132").Value = -5948.7245
$ws.Range("H137").Value = 50211.176
$ws.Range("I137").Value = 70255.375
$ws.Range("J137").Value = 4395.857
$ws.Range("K137").Value = 210766.125
$ws.Range("L137").Value = 13187.571
$ws.Range("M137").Value = -208216.125
$ws.Range("N137").Value = -18287.571
$ws.Range("H138").Value = 3385.923
$ws.Range("J138").Value = 3299.16
$ws.Range("L138").Value = 9897.48
$ws.Range("N138").Value = -20177.48
$ws.Range("H141").Value = 1522.2858
$ws.Range("I141").Value = 1523.4
$ws.Range("J141").Value = 1500
$ws.Range("K141").Value = 4570.200000000001
$ws.Range("L141").Value = 4500
$ws.Range("M141").Value = 609.7999999999993
$ws.Range("N141").Value = -14860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 11515.1
$ws.Range("I2").Value = 13629.75
$ws.Range("J2").Value = 3056.5
$ws.Range("K2").Value = 13629.75
$ws.Range("L2").Value = 3056.5
$ws.Range("M2").Value = -13516.75
$ws.Range("N2").Value = -3282.5
$ws.Range("H32").Value = 3992.468
$ws.Range("I32").Value = 3579.1
$ws.Range("K32").Value = 3579.1
$ws.Range("M32").Value = -3292.1
$ws.Range("H45").Value = 11167.75
$ws.Range("I45").Value = 22340.2
$ws.Range("J45").Value = 3187.4285
$ws.Range("K45").Value = 22340.2
$ws.Range("L45").Value = 3187.4285
$ws.Range("M45").Value = -21963.2
$ws.Range("N45").Value = -3941.4285
$ws.Range("H60").Value = 16050.5
$ws.Range("I60").Value = 16050.5
$ws.Range("K60").Value = 16050.5
$ws.Range("M60").Value = -15317.5
$ws.Range("H61").Value = 2641.0962
$ws.Range("I61").Value = 1816.7561
$ws.Range("J61").Value = 5713.636
$ws.Range("K61").Value = 1816.7561
$ws.Range("L61").Value = 5713.636
$ws.Range("M61").Value = -1604.7561
$ws.Range("N61").Value = -6137.636
$ws.Range("H63").Value = 3495.9092
$ws.Range("I63").Value = 3545.6
$ws.Range("K63").Value = 3545.6
$ws.Range("M63").Value = -2859.6
$ws.Range("H66").Value = 3495.9092
$ws.Range("I66").Value = 3545.6
$ws.Range("K66").Value = 17728
$ws.Range("M66").Value = -14296
$ws.Range("H74").Value = 106225.16
$ws.Range("I74").Value = 124515.69
$ws.Range("J74").Value = 8675.667
$ws.Range("K74").Value = 124515.69
$ws.Range("L74").Value = 8675.667
$ws.Range("M74").Value = -123641.69
$ws.Range("N74").Value = -10423.667
$ws.Range("H77").Value = 106225.16
$ws.Range("I77").Value = 124515.69
$ws.Range("J77").Value = 8675.667
$ws.Range("K77").Value = 622578.45
$ws.Range("L77").Value = 43378.335
$ws.Range("M77").Value = -618210.45
$ws.Range("N77").Value = -52114.335
$ws.Range("H88").Value = 2761.611
$ws.Range("I88").Value = 2594.25
$ws.Range("K88").Value = 2594.25
$ws.Range("M88").Value = -2188.25
$ws.Range("H91").Value = 2761.611
$ws.Range("I91").Value = 2594.25
$ws.Range("K91").Value = 2594.25
$ws.Range("M91").Value = -1190.25
$ws.Range("H97").Value = 1707.1154
$ws.Range("I97").Value = 1194.8889
$ws.Range("K97").Value = 1194.8889
$ws.Range("M97").Value = -698.8888999999999
$ws.Range("H102").Value = 2579.0833
$ws.Range("I102").Value = 2105.7778
$ws.Range("K102").Value = 2105.7778
$ws.Range("M102").Value = -483.7777999999998
$ws.Range("H110").Value = 5029.273
$ws.Range("I110").Value = 4849.9414
$ws.Range("J110").Value = 5639
$ws.Range("K110").Value = 4849.9414
$ws.Range("L110").Value = 5639
$ws.Range("M110").Value = -2804.9414
$ws.Range("N110").Value = -9729
$ws.Range("H116").Value = 11515.1
$ws.Range("I116").Value = 13629.75
$ws.Range("J116").Value = 3056.5
$ws.Range("K116").Value = 13629.75
$ws.Range("L116").Value = 3056.5
$ws.Range("M116").Value = -11335.75
$ws.Range("N116").Value = -7644.5
$ws.Range("H122").Value = 6588.4062
$ws.Range("I122").Value = 6406.3
$ws.Range("K122").Value = 19218.9
$ws.Range("M122").Value = -16768.9
$ws.Range("H126").Value = 8499.667
$ws.Range("I126").Value = 8499.667
$ws.Range("K126").Value = 25499.001
$ws.Range("M126").Value = -23029.001
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").Value = $null
$ws.Range("H132").Value = 3809.5625
$ws.Range("I132").Value = 3575.8928
$ws.Range("J132").Value = 5445.25
$ws.Range("K132").Value = 10727.6784
$ws.Range("L132").Value = 16335.75
$ws.Range("M132").Value = -8197.6784
$ws.Range("N132").Value = -21395.75
$ws.Range("H136").Value = 2641.0962
$ws.Range("I136").Value = 1816.7561
$ws.Range("J136").Value = 5713.636
$ws.Range("K136").Value = 5450.2683
$ws.Range("L136").Value = 17140.908
$ws.Range("M136").Value = -2900.2683
$ws.Range("N136").Value = -22240.908

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 11515.1
$ws.Range("I3").Value = 13629.75
$ws.Range("J3").Value = 3056.5
$ws.Range("K3").Value = 13629.75
$ws.Range("L3").Value = 3056.5
$ws.Range("M3").Value = -13515.75
$ws.Range("N3").Value = -3284.5
$ws.Range("H12").Value = 1030.6666
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 1030.6666
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 1030.6666
$ws.Range("M12").Value = $null
$ws.Range("N12").Value = -1366.6666
$ws.Range("H18").Value = 1500
$ws.Range("J18").Value = 1500
$ws.Range("L18").Value = 1500
$ws.Range("N18").Value = -2558
$ws.Range("H20").Value = 3746.8462
$ws.Range("I20").Value = 3526.95
$ws.Range("J20").Value = 4479.8335
$ws.Range("K20").Value = 3526.95
$ws.Range("L20").Value = 4479.8335
$ws.Range("M20").Value = -3279.95
$ws.Range("N20").Value = -4973.8335
$ws.Range("H22").Value = 496.33334
$ws.Range("I22").Value = 521.0714
$ws.Range("J22").Value = 150
$ws.Range("K22").Value = 521.0714
$ws.Range("L22").Value = 150
$ws.Range("M22").Value = -348.0714
$ws.Range("N22").Value = -496
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = $null
$ws.Range("N37").Value = $null
$ws.Range("H60").Value = 41495
$ws.Range("J60").Value = 41495
$ws.Range("L60").Value = 41495
$ws.Range("N60").Value = -42693
$ws.Range("H80").Value = 649.54285
$ws.Range("J80").Value = 318.15
$ws.Range("L80").Value = 318.15
$ws.Range("N80").Value = -2314.15
$ws.Range("H83").Value = 649.54285
$ws.Range("J83").Value = 318.15
$ws.Range("L83").Value = 1590.75
$ws.Range("N83").Value = -11574.75
$ws.Range("H86").Value = 2573.578
$ws.Range("I86").Value = 2233.4644
$ws.Range("J86").Value = 3133.7646
$ws.Range("K86").Value = 2233.4644
$ws.Range("L86").Value = 3133.7646
$ws.Range("M86").Value = -1110.4644
$ws.Range("N86").Value = -5379.7646
$ws.Range("H89").Value = 2573.578
$ws.Range("I89").Value = 2233.4644
$ws.Range("J89").Value = 3133.7646
$ws.Range("K89").Value = 11167.322
$ws.Range("L89").Value = 15668.823
$ws.Range("M89").Value = -5551.322
$ws.Range("N89").Value = -26900.823
$ws.Range("H94").Value = 687.75
$ws.Range("I94").Value = 609.8947
$ws.Range("K94").Value = 609.8947
$ws.Range("M94").Value = -158.8946999999999
$ws.Range("H99").Value = 5005
$ws.Range("I99").Value = 5037.5
$ws.Range("K99").Value = 5037.5
$ws.Range("M99").Value = -3539.5
$ws.Range("H105").Value = 2161.724
$ws.Range("I105").Value = 2043.9584
$ws.Range("K105").Value = 2043.9584
$ws.Range("M105").Value = -296.9584
$ws.Range("H128").Value = 59330.715
$ws.Range("I128").Value = 59330.715
$ws.Range("K128").Value = 177992.145
$ws.Range("M128").Value = -175502.145
$ws.Range("H134").Value = 4444.2
$ws.Range("I134").Value = 4483.5415
$ws.Range("K134").Value = 13450.6245
$ws.Range("M134").Value = -10915.6245

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1544.8334
$ws.Range("I16").Value = 788.4375
$ws.Range("K16").Value = 788.4375
$ws.Range("M16").Value = -501.4375
$ws.Range("H31").Value = 171703.88
$ws.Range("I31").Value = 372286.97
$ws.Range("J31").Value = 2461.9062
$ws.Range("K31").Value = 372286.97
$ws.Range("L31").Value = 2461.9062
$ws.Range("M31").Value = -371991.97
$ws.Range("N31").Value = -3051.9062
$ws.Range("H34").Value = 171703.88
$ws.Range("I34").Value = 372286.97
$ws.Range("J34").Value = 2461.9062
$ws.Range("K34").Value = 372286.97
$ws.Range("L34").Value = 2461.9062
$ws.Range("M34").Value = -372084.97
$ws.Range("N34").Value = -2865.9062
$ws.Range("H58").Value = 2324.5
$ws.Range("I58").Value = 2288.6667
$ws.Range("J58").Value = 2378.25
$ws.Range("K58").Value = 2288.6667
$ws.Range("L58").Value = 2378.25
$ws.Range("M58").Value = -2085.6667
$ws.Range("N58").Value = -2784.25
$ws.Range("H69").Value = 19841
$ws.Range("I69").Value = 19841
$ws.Range("K69").Value = 19841
$ws.Range("M69").Value = -19092
$ws.Range("H72").Value = 19841
$ws.Range("I72").Value = 19841
$ws.Range("K72").Value = 59523
$ws.Range("M72").Value = -55779
$ws.Range("H80").Value = 30000
$ws.Range("J80").Value = 30000
$ws.Range("L80").Value = 30000
$ws.Range("N80").Value = -32246
$ws.Range("H83").Value = 30000
$ws.Range("J83").Value = 30000
$ws.Range("L83").Value = 90000
$ws.Range("N83").Value = -101232
$ws.Range("H94").Value = 1461.6666
$ws.Range("J94").Value = 1477.7693
$ws.Range("L94").Value = 1477.7693
$ws.Range("N94").Value = -2379.7693
$ws.Range("H99").Value = 380962.62
$ws.Range("I99").Value = 914596.4
$ws.Range("J99").Value = 14089.4375
$ws.Range("K99").Value = 914596.4
$ws.Range("L99").Value = 14089.4375
$ws.Range("M99").Value = -913098.4
$ws.Range("N99").Value = -17085.4375
$ws.Range("H113").Value = 1544.8334
$ws.Range("I113").Value = 788.4375
$ws.Range("K113").Value = 788.4375
$ws.Range("M113").Value = 1381.5625
$ws.Range("H126").Value = 380962.62
$ws.Range("I126").Value = 914596.4
$ws.Range("J126").Value = 14089.4375
$ws.Range("K126").Value = 2743789.2
$ws.Range("L126").Value = 42268.3125
$ws.Range("M126").Value = -2741319.2
$ws.Range("N126").Value = -47208.3125
$ws.Range("H132").Value = 4117.0527
$ws.Range("I132").Value = 2248.7646
$ws.Range("K132").Value = 6746.293799999999
$ws.Range("M132").Value = -4216.293799999999
$ws.Range("H134").Value = 3040.037
$ws.Range("I134").Value = 2542.8262
$ws.Range("K134").Value = 7628.4786
$ws.Range("M134").Value = -5093.4786
$ws.Range("H136").Value = 2324.5
$ws.Range("I136").Value = 2288.6667
$ws.Range("J136").Value = 2378.25
$ws.Range("K136").Value = 6866.000100000001
$ws.Range("L136").Value = 7134.75
$ws.Range("M136").Value = -4316.000100000001
$ws.Range("N136").Value = -12234.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1775217.1
$ws.Range("I4").Value = 1825813.6
$ws.Range("K4").Value = 5477440.800000001
$ws.Range("M4").Value = -5477328.800000001
$ws.Range("H12").Value = 3826.15
$ws.Range("I12").Value = 714.8333
$ws.Range("J12").Value = 5159.5713
$ws.Range("K12").Value = 2144.4999
$ws.Range("L12").Value = 15478.7139
$ws.Range("M12").Value = -1971.4999
$ws.Range("N12").Value = -15824.7139
$ws.Range("H17").Value = 184.57143
$ws.Range("J17").Value = 215.66667
$ws.Range("L17").Value = 647.00001
$ws.Range("N17").Value = -985.00001
$ws.Range("H26").Value = 3874.25
$ws.Range("J26").Value = 5121.5
$ws.Range("L26").Value = 15364.5
$ws.Range("N26").Value = -15940.5
$ws.Range("H33").Value = 2131.125
$ws.Range("I33").Value = 3902.8462
$ws.Range("J33").Value = 37.272728
$ws.Range("K33").Value = 23417.0772
$ws.Range("L33").Value = 223.636368
$ws.Range("M33").Value = -23134.0772
$ws.Range("N33").Value = -789.636368
$ws.Range("H34").Value = 960.53125
$ws.Range("I34").Value = 629.8077
$ws.Range("J34").Value = 2393.6667
$ws.Range("K34").Value = 1889.4231
$ws.Range("L34").Value = 7181.000100000001
$ws.Range("M34").Value = -1805.4231
$ws.Range("N34").Value = -7349.000100000001
$ws.Range("H35").Value = 2518
$ws.Range("I35").Value = 1363.3334
$ws.Range("J35").Value = 4250
$ws.Range("K35").Value = 4090.0002
$ws.Range("L35").Value = 12750
$ws.Range("M35").Value = -3802.0002
$ws.Range("N35").Value = -13326
$ws.Range("H37").Value = 76985480
$ws.Range("J37").Value = 76985480
$ws.Range("L37").Value = 230956440
$ws.Range("N37").Value = -230956664
$ws.Range("H39").Value = 811.5
$ws.Range("I39").Value = 696.3684
$ws.Range("K39").Value = 2089.1052
$ws.Range("M39").Value = -1795.1052
$ws.Range("H44").Value = 2146.0625
$ws.Range("J44").Value = 2934.0908
$ws.Range("L44").Value = 8802.2724
$ws.Range("N44").Value = -9598.2724
$ws.Range("H46").Value = 6112349.5
$ws.Range("I46").Value = 1098.75
$ws.Range("K46").Value = 3296.25
$ws.Range("M46").Value = -3205.25
$ws.Range("H55").Value = 6711.8887
$ws.Range("I55").Value = 2999
$ws.Range("J55").Value = 7772.7144
$ws.Range("K55").Value = 8997
$ws.Range("L55").Value = 23318.1432
$ws.Range("M55").Value = -8820
$ws.Range("N55").Value = -23672.1432
$ws.Range("H63").Value = 3000
$ws.Range("J63").Value = 5000
$ws.Range("L63").Value = 15000
$ws.Range("N63").Value = -16498
$ws.Range("H66").Value = 3000
$ws.Range("J66").Value = 5000
$ws.Range("L66").Value = 45000
$ws.Range("N66").Value = -52488
$ws.Range("H68").Value = 4388023
$ws.Range("J68").Value = 2479.7144
$ws.Range("L68").Value = 7439.1432
$ws.Range("N68").Value = -9061.143199999999
$ws.Range("H71").Value = 4388023
$ws.Range("J71").Value = 2479.7144
$ws.Range("L71").Value = 22317.4296
$ws.Range("N71").Value = -30429.4296
$ws.Range("H99").Value = 13028.77
$ws.Range("I99").Value = 14924.875
$ws.Range("J99").Value = 9995
$ws.Range("K99").Value = 44774.625
$ws.Range("L99").Value = 29985
$ws.Range("M99").Value = -42528.625
$ws.Range("N99").Value = -34477
$ws.Range("H107").Value = 925.3143
$ws.Range("I107").Value = 355.55
$ws.Range("K107").Value = 1066.65
$ws.Range("M107").Value = 853.3499999999999
$ws.Range("H113").Value = 389.39285
$ws.Range("I113").Value = 202.33333
$ws.Range("K113").Value = 606.99999
$ws.Range("M113").Value = 1563.00001
$ws.Range("H116").Value = 1734.5
$ws.Range("I116").Value = 1500
$ws.Range("J116").Value = 1969
$ws.Range("K116").Value = 4500
$ws.Range("L116").Value = 5907
$ws.Range("M116").Value = -1058
$ws.Range("N116").Value = -12791
$ws.Range("H122").Value = 794.1667
$ws.Range("J122").Value = 816.6
$ws.Range("L122").Value = 7349.400000000001
$ws.Range("N122").Value = -12249.4
$ws.Range("H123").Value = 4312.4165
$ws.Range("I123").Value = 4179.9
$ws.Range("J123").Value = 4975
$ws.Range("K123").Value = 12539.7
$ws.Range("L123").Value = 14925
$ws.Range("M123").Value = -10089.7
$ws.Range("N123").Value = -19825
$ws.Range("H131").Value = 13515107
$ws.Range("I131").Value = 38462700
$ws.Range("J131").Value = 1828.0834
$ws.Range("K131").Value = 115388100
$ws.Range("L131").Value = 5484.2502
$ws.Range("M131").Value = -115383060
$ws.Range("N131").Value = -15564.2502
$ws.Range("H132").Value = 5797.3794
$ws.Range("I132").Value = 8601.1875
$ws.Range("J132").Value = 2346.5386
$ws.Range("K132").Value = 77410.6875
$ws.Range("L132").Value = 21118.8474
$ws.Range("M132").Value = -74880.6875
$ws.Range("N132").Value = -26178.8474
$ws.Range("H140").Value = 937.619
$ws.Range("I140").Value = 809.5
$ws.Range("J140").Value = 3500
$ws.Range("K140").Value = 2428.5
$ws.Range("L140").Value = 10500
$ws.Range("M140").Value = 2751.5
$ws.Range("N140").Value = -20860
$ws.Range("H141").Value = 3375.5
$ws.Range("I141").Value = 2325.7
$ws.Range("K141").Value = 6977.099999999999
$ws.Range("M141").Value = -1797.099999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 20000000
$ws.Range("I21").Value = 20000000
$ws.Range("K21").Value = 20000000
$ws.Range("M21").Value = -19999827
$ws.Range("H30").Value = 20000000
$ws.Range("I30").Value = 20000000
$ws.Range("K30").Value = 20000000
$ws.Range("M30").Value = -19999895
$ws.Range("H36").Value = 10000
$ws.Range("I36").Value = 10000
$ws.Range("K36").Value = 10000
$ws.Range("M36").Value = -9515
$ws.Range("H44").Value = 20000
$ws.Range("I44").Value = 20000
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 20000
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -19404
$ws.Range("N44").Value = $null
$ws.Range("H58").Value = 22203
$ws.Range("I58").Value = 20347
$ws.Range("J58").Value = 24987
$ws.Range("K58").Value = 20347
$ws.Range("L58").Value = 24987
$ws.Range("M58").Value = -20070
$ws.Range("N58").Value = -25541
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = $null
$ws.Range("N70").Value = $null
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = $null
$ws.Range("N73").Value = $null
$ws.Range("H92").Value = 10875.5
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 10875.5
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 10875.5
$ws.Range("M92").Value = $null
$ws.Range("N92").Value = -14619.5
$ws.Range("H97").Value = 1802.3158
$ws.Range("I97").Value = 1634.0625
$ws.Range("K97").Value = 1634.0625
$ws.Range("M97").Value = -1138.0625
$ws.Range("H102").Value = 34581.906
$ws.Range("I102").Value = 3003
$ws.Range("J102").Value = 129318.625
$ws.Range("K102").Value = 3003
$ws.Range("L102").Value = 129318.625
$ws.Range("M102").Value = -1381
$ws.Range("N102").Value = -132562.625
$ws.Range("H126").Value = 8339.058
$ws.Range("I126").Value = 8646.161
$ws.Range("J126").Value = 5959
$ws.Range("K126").Value = 25938.483
$ws.Range("L126").Value = 17877
$ws.Range("M126").Value = -23468.483
$ws.Range("N126").Value = -22817
$ws.Range("H131").Value = 49995
$ws.Range("J131").Value = 49995
$ws.Range("L131").Value = 49995
$ws.Range("N131").Value = -60075
$ws.Range("H132").Value = 53843.8
$ws.Range("I132").Value = 70387.07
$ws.Range("J132").Value = 4214
$ws.Range("K132").Value = 211161.21
$ws.Range("L132").Value = 12642
$ws.Range("M132").Value = -208631.21
$ws.Range("N132").Value = -17702

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5036.25
$ws.Range("I7").Value = 4066.2222
$ws.Range("K7").Value = 4066.2222
$ws.Range("M7").Value = -3954.2222
$ws.Range("H16").Value = 31251234
$ws.Range("I16").Value = 35715572
$ws.Range("J16").Value = 881
$ws.Range("K16").Value = 35715572
$ws.Range("L16").Value = 881
$ws.Range("M16").Value = -35715402
$ws.Range("N16").Value = -1221
$ws.Range("H17").Value = 3834.8333
$ws.Range("I17").Value = 3377.25
$ws.Range("J17").Value = 4750
$ws.Range("K17").Value = 3377.25
$ws.Range("L17").Value = 4750
$ws.Range("M17").Value = -3207.25
$ws.Range("N17").Value = -5090
$ws.Range("H22").Value = 964
$ws.Range("I22").Value = 950
$ws.Range("K22").Value = 950
$ws.Range("M22").Value = -655
$ws.Range("H27").Value = 964
$ws.Range("I27").Value = 950
$ws.Range("K27").Value = 950
$ws.Range("M27").Value = -843
$ws.Range("H33").Value = 10017
$ws.Range("J33").Value = 10017
$ws.Range("L33").Value = 10017
$ws.Range("N33").Value = -10597
$ws.Range("H40").Value = 4908.8237
$ws.Range("I40").Value = 4123.5
$ws.Range("K40").Value = 4123.5
$ws.Range("M40").Value = -3987.5
$ws.Range("H46").Value = 4108.154
$ws.Range("I46").Value = 3142.6
$ws.Range("J46").Value = 7326.6665
$ws.Range("K46").Value = 3142.6
$ws.Range("L46").Value = 7326.6665
$ws.Range("M46").Value = -2954.6
$ws.Range("N46").Value = -7702.6665
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").Value = $null
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").Value = $null
$ws.Range("H55").Value = 6048.5356
$ws.Range("I55").Value = 716.8
$ws.Range("J55").Value = 19377.875
$ws.Range("K55").Value = 716.8
$ws.Range("L55").Value = 19377.875
$ws.Range("M55").Value = -543.8
$ws.Range("N55").Value = -19723.875
$ws.Range("H93").Value = 3116
$ws.Range("I93").Value = 1852.1538
$ws.Range("K93").Value = 1852.1538
$ws.Range("M93").Value = -604.1538
$ws.Range("H100").Value = 25923.25
$ws.Range("I100").Value = 34333.332
$ws.Range("J100").Value = 693
$ws.Range("K100").Value = 34333.332
$ws.Range("L100").Value = 693
$ws.Range("M100").Value = -33792.332
$ws.Range("N100").Value = -1775
$ws.Range("H122").Value = 5163.1665
$ws.Range("I122").Value = 4247.75
$ws.Range("J122").Value = 6994
$ws.Range("K122").Value = 12743.25
$ws.Range("L122").Value = 20982
$ws.Range("M122").Value = -10293.25
$ws.Range("N122").Value = -25882
$ws.Range("H126").Value = 5036.25
$ws.Range("I126").Value = 4066.2222
$ws.Range("K126").Value = 12198.6666
$ws.Range("M126").Value = -9728.6666
$ws.Range("H132").Value = 3019.8333
$ws.Range("I132").Value = 2728.0881
$ws.Range("J132").Value = 4259.75
$ws.Range("K132").Value = 8184.2643
$ws.Range("L132").Value = 12779.25
$ws.Range("M132").Value = -5654.2643
$ws.Range("N132").Value = -17839.25
$ws.Range("H136").Value = 13360.5
$ws.Range("I136").Value = 22222
$ws.Range("J136").Value = 4499
$ws.Range("K136").Value = 66666
$ws.Range("L136").Value = 13497
$ws.Range("M136").Value = -64116
$ws.Range("N136").Value = -18597

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").Value = $null
$ws.Range("H46").Value = 94452.75
$ws.Range("J46").Value = 94452.75
$ws.Range("L46").Value = 94452.75
$ws.Range("N46").Value = -94914.75
$ws.Range("H51").Value = 11998.6
$ws.Range("I51").Value = 11998.6
$ws.Range("K51").Value = 11998.6
$ws.Range("M51").Value = -11488.6
$ws.Range("H52").Value = 23999
$ws.Range("I52").Value = 23999
$ws.Range("K52").Value = 23999
$ws.Range("M52").Value = -23773
$ws.Range("H54").Value = 19990
$ws.Range("J54").Value = 19990
$ws.Range("L54").Value = 19990
$ws.Range("N54").Value = -21030
$ws.Range("H62").Value = 27787054
$ws.Range("I62").Value = 8125
$ws.Range("J62").Value = 50010196
$ws.Range("K62").Value = 8125
$ws.Range("L62").Value = 50010196
$ws.Range("M62").Value = -7501
$ws.Range("N62").Value = -50011444
$ws.Range("H65").Value = 27787054
$ws.Range("I65").Value = 8125
$ws.Range("J65").Value = 50010196
$ws.Range("K65").Value = 40625
$ws.Range("L65").Value = 250050980
$ws.Range("M65").Value = -37505
$ws.Range("N65").Value = -250057220
$ws.Range("H81").Value = 7666.3335
$ws.Range("I81").Value = 7666.3335
$ws.Range("K81").Value = 15332.667
$ws.Range("M81").Value = -14271.667
$ws.Range("H84").Value = 7666.3335
$ws.Range("I84").Value = 7666.3335
$ws.Range("K84").Value = 76663.33499999999
$ws.Range("M84").Value = -71359.33499999999
$ws.Range("H97").Value = 21499.5
$ws.Range("J97").Value = 21499.5
$ws.Range("L97").Value = 21499.5
$ws.Range("N97").Value = -23481.5
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = $null
$ws.Range("N100").Value = $null
$ws.Range("H107").Value = 512.4706
$ws.Range("I107").Value = 486
$ws.Range("J107").Value = 561
$ws.Range("K107").Value = 1458
$ws.Range("L107").Value = 1683
$ws.Range("M107").Value = 462
$ws.Range("N107").Value = -5523
$ws.Range("H113").Value = 659.4091
$ws.Range("I113").Value = 627.7895
$ws.Range("K113").Value = 1883.3685
$ws.Range("M113").Value = 286.6315
$ws.Range("H122").Value = 3722.8064
$ws.Range("I122").Value = 3642
$ws.Range("K122").Value = 10926
$ws.Range("M122").Value = -8476
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = $null
$ws.Range("H126").Value = 2906.7585
$ws.Range("I126").Value = 2665.6956
$ws.Range("J126").Value = 3830.8333
$ws.Range("K126").Value = 7997.0868
$ws.Range("L126").Value = 11492.4999
$ws.Range("M126").Value = -5527.0868
$ws.Range("N126").Value = -16432.4999
$ws.Range("H132").Value = 1953.7805
$ws.Range("I132").Value = 1887.3077
$ws.Range("J132").Value = 3250
$ws.Range("K132").Value = 5661.9231
$ws.Range("L132").Value = 9750
$ws.Range("M132").Value = -3131.9231
$ws.Range("N132").Value = -14810
$ws.Range("H134").Value = 94452.75
$ws.Range("J134").Value = 94452.75
$ws.Range("L134").Value = 283358.25
$ws.Range("N134").Value = -288428.25
$ws.Range("H136").Value = 240427.33
$ws.Range("I136").Value = 246169.47
$ws.Range("K136").Value = 738508.41
$ws.Range("M136").Value = -735958.41
